$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 279 (pushing existing rows 279+ down to 281+)
$ws.Rows("279:280").Insert()

# Populate the two newly inserted rows with their data
$ws.Range("A279").Value = 5
$ws.Range("B279").Value = 'Macroferia Regional de Talca'
$ws.Range("C279").Value = 'Maule'
$ws.Range("D279").Value = 44524
$ws.Range("E279").Value = 7
$ws.Range("F279").Value = 'Fruta'
$ws.Range("G279").Value = 100102
$ws.Range("H279").Value = 'Cítricos'
$ws.Range("I279").Value = 100102005
$ws.Range("J279").Value = 'Naranja'
$ws.Range("K279").Value = 'Lane Late'
$ws.Range("L279").Value = 'Primera'
$ws.Range("M279").Value = 250
$ws.Range("N279").Value = 8000
$ws.Range("O279").Value = 8000
$ws.Range("P279").Value = 8000
$ws.Range("Q279").Value = '$/bandeja 15 kilos granel'
$ws.Range("R279").Value = 'Región de O''Higgins'
$ws.Range("S279").Value = 533
$ws.Range("T279").Value = 15
$ws.Range("A280").Value = 5
$ws.Range("B280").Value = 'Macroferia Regional de Talca'
$ws.Range("C280").Value = 'Maule'
$ws.Range("D280").Value = 44524
$ws.Range("E280").Value = 7
$ws.Range("F280").Value = 'Fruta'
$ws.Range("G280").Value = 100102
$ws.Range("H280").Value = 'Cítricos'
$ws.Range("I280").Value = 100102005
$ws.Range("J280").Value = 'Naranja'
$ws.Range("K280").Value = 'Navel Late'
$ws.Range("L280").Value = 'Primera'
$ws.Range("M280").Value = 360
$ws.Range("N280").Value = 8000
$ws.Range("O280").Value = 8000
$ws.Range("P280").Value = 8000
$ws.Range("Q280").Value = '$/bandeja 15 kilos granel'
$ws.Range("R280").Value = 'Región de O''Higgins'
$ws.Range("S280").Value = 533
$ws.Range("T280").Value = 15
